# Change order of scenarios (Negative interventions at the end)

$wb = $excel.ActiveWorkbook
$wsConstants = $wb.Worksheets.Item("constants")
$wsTimeVariants = $wb.Worksheets.Item("time_variants")

$xlPasteFormats = -4122
$xlNone = -4142

# ---------------------------------------------------------------------------
# constants sheet
# ---------------------------------------------------------------------------

# tb_n_contact value changed; old value moved to a new column E2
$wsConstants.Range("D2").Copy()
$wsConstants.Range("E2").PasteSpecial($xlPasteFormats)
$wsConstants.Range("E2").Value = 2.7122322346359198
$wsConstants.Range("B2").Value = 2.2999999999999998

# start_time changed from 1990 to 2010
$wsConstants.Range("B17").Value = 2010

# ---------------------------------------------------------------------------
# time_variants sheet -- re-order scenario intervention columns so that the
# negative intervention scenario ends up in the last populated columns
# ---------------------------------------------------------------------------

# Row 5
$wsTimeVariants.Range("S5").Value = 100
$wsTimeVariants.Range("U5").Copy()
$wsTimeVariants.Range("V5").PasteSpecial($xlPasteFormats)
$wsTimeVariants.Range("V5").ClearContents()

# Row 8
$wsTimeVariants.Range("S8").Copy()
$wsTimeVariants.Range("T8").PasteSpecial($xlPasteFormats)
$wsTimeVariants.Range("T8").ClearContents()
$wsTimeVariants.Range("S8").ClearContents()
$wsTimeVariants.Range("V8").Value = 0

# Row 9
$wsTimeVariants.Range("T9").Copy()
$wsTimeVariants.Range("V9").PasteSpecial($xlPasteFormats)
$wsTimeVariants.Range("V9").Value = 0
$wsTimeVariants.Range("S9").Copy()
$wsTimeVariants.Range("U9").PasteSpecial($xlPasteFormats)
$wsTimeVariants.Range("U9").Value = 0
$wsTimeVariants.Range("S9").ClearContents()
$wsTimeVariants.Range("T9").ClearContents()

# Row 10
$wsTimeVariants.Range("U10").Copy()
$wsTimeVariants.Range("T10").PasteSpecial($xlPasteFormats)
$wsTimeVariants.Range("T10").Value = 100
$wsTimeVariants.Range("U10").Copy()
$wsTimeVariants.Range("U10").PasteSpecial($xlPasteFormats)
$wsTimeVariants.Range("U10").ClearContents()
$wsTimeVariants.Range("S10").Copy()
$wsTimeVariants.Range("U10").PasteSpecial($xlPasteFormats)

# Row 11
$wsTimeVariants.Range("U11").Copy()
$wsTimeVariants.Range("T11").PasteSpecial($xlPasteFormats)
$wsTimeVariants.Range("T11").Value = 100
$wsTimeVariants.Range("S11").Copy()
$wsTimeVariants.Range("U11").PasteSpecial($xlPasteFormats)

Write-Output "done"
